$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the now-unused trailing rows (32:33) so the sheet dimension becomes A1:Q31
$ws.Range("A32:Q33").Delete()

# Rewrite the negative/positive word-stat tables with the recomputed (min-count=5) values
$ws.Range("A1").Value = "negative"
$ws.Range("J1").Value = "positive"
$ws.Range("A2").Value = "name"
$ws.Range("B2").Value = "anchor score"
$ws.Range("C2").Value = "type occurences"
$ws.Range("D2").Value = "total occurences"
$ws.Range("E2").Value = "+%"
$ws.Range("F2").Value = "-%"
$ws.Range("G2").Value = "both"
$ws.Range("H2").Value = "normal"
$ws.Range("J2").Value = "name"
$ws.Range("K2").Value = "anchor score"
$ws.Range("L2").Value = "type occurences"
$ws.Range("M2").Value = "total occurences"
$ws.Range("N2").Value = "+%"
$ws.Range("O2").Value = "-%"
$ws.Range("P2").Value = "both"
$ws.Range("Q2").Value = "normal"
$ws.Range("A3").Value = "poorly"
$ws.Range("B3").Value = 0.9565217391304348
$ws.Range("C3").Value = 44
$ws.Range("D3").Value = 44
$ws.Range("E3").Value = 0
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = $false
$ws.Range("H3").Value = 2
$ws.Range("J3").Value = "wonderful"
$ws.Range("K3").Value = 0.8928571428571429
$ws.Range("L3").Value = 50
$ws.Range("M3").Value = 50
$ws.Range("N3").Value = 1
$ws.Range("O3").Value = 0
$ws.Range("P3").Value = $false
$ws.Range("Q3").Value = 6
$ws.Range("A4").Value = "disappointing"
$ws.Range("B4").Value = 0.7954545454545454
$ws.Range("C4").Value = 35
$ws.Range("D4").Value = 35
$ws.Range("E4").Value = 0
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = $false
$ws.Range("H4").Value = 9
$ws.Range("J4").Value = "awesome"
$ws.Range("K4").Value = 0.8615384615384616
$ws.Range("L4").Value = 56
$ws.Range("M4").Value = 56
$ws.Range("N4").Value = 1
$ws.Range("O4").Value = 0
$ws.Range("P4").Value = $false
$ws.Range("Q4").Value = 9
$ws.Range("A5").Value = "poor"
$ws.Range("B5").Value = 0.7183098591549296
$ws.Range("C5").Value = 51
$ws.Range("D5").Value = 51
$ws.Range("E5").Value = 0
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = $false
$ws.Range("H5").Value = 20
$ws.Range("J5").Value = "favorite"
$ws.Range("K5").Value = 0.8279569892473119
$ws.Range("L5").Value = 77
$ws.Range("M5").Value = 77
$ws.Range("N5").Value = 1
$ws.Range("O5").Value = 0
$ws.Range("P5").Value = $false
$ws.Range("Q5").Value = 16
$ws.Range("A6").Value = "disappointed"
$ws.Range("B6").Value = 0.7096774193548387
$ws.Range("C6").Value = 132
$ws.Range("D6").Value = 132
$ws.Range("E6").Value = 0
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = $false
$ws.Range("H6").Value = 54
$ws.Range("J6").Value = "excellent"
$ws.Range("K6").Value = 0.796875
$ws.Range("L6").Value = 51
$ws.Range("M6").Value = 51
$ws.Range("N6").Value = 1
$ws.Range("O6").Value = 0
$ws.Range("P6").Value = $false
$ws.Range("Q6").Value = 13
$ws.Range("A7").Value = "however"
$ws.Range("B7").Value = 0.671875
$ws.Range("C7").Value = 43
$ws.Range("D7").Value = 43
$ws.Range("E7").Value = 0
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = $false
$ws.Range("H7").Value = 21
$ws.Range("J7").Value = "classic"
$ws.Range("K7").Value = 0.660377358490566
$ws.Range("L7").Value = 35
$ws.Range("M7").Value = 35
$ws.Range("N7").Value = 1
$ws.Range("O7").Value = 0
$ws.Range("P7").Value = $false
$ws.Range("Q7").Value = 18
$ws.Range("A8").Value = "broke"
$ws.Range("B8").Value = 0.6310679611650486
$ws.Range("C8").Value = 130
$ws.Range("D8").Value = 130
$ws.Range("E8").Value = 0
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = $false
$ws.Range("H8").Value = 76
$ws.Range("J8").Value = "thank"
$ws.Range("K8").Value = 0.5797101449275363
$ws.Range("L8").Value = 40
$ws.Range("M8").Value = 40
$ws.Range("N8").Value = 1
$ws.Range("O8").Value = 0
$ws.Range("P8").Value = $false
$ws.Range("Q8").Value = 29
$ws.Range("A9").Value = "waste"
$ws.Range("B9").Value = 0.6148648648648649
$ws.Range("C9").Value = 91
$ws.Range("D9").Value = 91
$ws.Range("E9").Value = 0
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = $false
$ws.Range("H9").Value = 57
$ws.Range("J9").Value = "love"
$ws.Range("K9").Value = 0.5265423242467718
$ws.Range("L9").Value = 367
$ws.Range("M9").Value = 367
$ws.Range("N9").Value = 1
$ws.Range("O9").Value = 0
$ws.Range("P9").Value = $false
$ws.Range("Q9").Value = 330
$ws.Range("A10").Value = "junk"
$ws.Range("B10").Value = 0.6
$ws.Range("C10").Value = 33
$ws.Range("D10").Value = 33
$ws.Range("E10").Value = 0
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = $false
$ws.Range("H10").Value = 22
$ws.Range("J10").Value = "loves"
$ws.Range("K10").Value = 0.4979253112033195
$ws.Range("L10").Value = 240
$ws.Range("M10").Value = 240
$ws.Range("N10").Value = 1
$ws.Range("O10").Value = 0
$ws.Range("P10").Value = $false
$ws.Range("Q10").Value = 242
$ws.Range("A11").Value = "smaller"
$ws.Range("B11").Value = 0.5882352941176471
$ws.Range("C11").Value = 70
$ws.Range("D11").Value = 70
$ws.Range("E11").Value = 0
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = $false
$ws.Range("H11").Value = 49
$ws.Range("J11").Value = "great"
$ws.Range("K11").Value = 0.4540983606557377
$ws.Range("L11").Value = 554
$ws.Range("M11").Value = 554
$ws.Range("N11").Value = 1
$ws.Range("O11").Value = 0
$ws.Range("P11").Value = $false
$ws.Range("Q11").Value = 666
$ws.Range("A12").Value = "small"
$ws.Range("B12").Value = 0.4927536231884058
$ws.Range("C12").Value = 170
$ws.Range("D12").Value = 170
$ws.Range("E12").Value = 0
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = $false
$ws.Range("H12").Value = 175
$ws.Range("J12").Value = "perfect"
$ws.Range("K12").Value = 0.3433734939759036
$ws.Range("L12").Value = 57
$ws.Range("M12").Value = 57
$ws.Range("N12").Value = 1
$ws.Range("O12").Value = 0
$ws.Range("P12").Value = $false
$ws.Range("Q12").Value = 109
$ws.Range("A13").Value = "apart"
$ws.Range("B13").Value = 0.4315789473684211
$ws.Range("C13").Value = 41
$ws.Range("D13").Value = 41
$ws.Range("E13").Value = 0
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = $false
$ws.Range("H13").Value = 54
$ws.Range("J13").Value = "loved"
$ws.Range("K13").Value = 0.3425076452599388
$ws.Range("L13").Value = 112
$ws.Range("M13").Value = 112
$ws.Range("N13").Value = 1
$ws.Range("O13").Value = 0
$ws.Range("P13").Value = $false
$ws.Range("Q13").Value = 215
$ws.Range("A14").Value = "broken"
$ws.Range("B14").Value = 0.3855421686746988
$ws.Range("C14").Value = 32
$ws.Range("D14").Value = 32
$ws.Range("E14").Value = 0
$ws.Range("F14").Value = 1
$ws.Range("G14").Value = $false
$ws.Range("H14").Value = 51
$ws.Range("J14").Value = "best"
$ws.Range("K14").Value = 0.3083333333333333
$ws.Range("L14").Value = 37
$ws.Range("M14").Value = 37
$ws.Range("N14").Value = 1
$ws.Range("O14").Value = 0
$ws.Range("P14").Value = $false
$ws.Range("Q14").Value = 83
$ws.Range("A15").Value = "cheap"
$ws.Range("B15").Value = 0.3791469194312796
$ws.Range("C15").Value = 80
$ws.Range("D15").Value = 80
$ws.Range("E15").Value = 0
$ws.Range("F15").Value = 1
$ws.Range("G15").Value = $false
$ws.Range("H15").Value = 131
$ws.Range("J15").Value = "friends"
$ws.Range("K15").Value = 0.2962962962962963
$ws.Range("L15").Value = 56
$ws.Range("M15").Value = 56
$ws.Range("N15").Value = 1
$ws.Range("O15").Value = 0
$ws.Range("P15").Value = $false
$ws.Range("Q15").Value = 133
$ws.Range("A16").Value = "plastic"
$ws.Range("B16").Value = 0.3622047244094488
$ws.Range("C16").Value = 46
$ws.Range("D16").Value = 46
$ws.Range("E16").Value = 0
$ws.Range("F16").Value = 1
$ws.Range("G16").Value = $false
$ws.Range("H16").Value = 81
$ws.Range("J16").Value = "christmas"
$ws.Range("K16").Value = 0.1807228915662651
$ws.Range("L16").Value = 45
$ws.Range("M16").Value = 45
$ws.Range("N16").Value = 1
$ws.Range("O16").Value = 0
$ws.Range("P16").Value = $false
$ws.Range("Q16").Value = 204
$ws.Range("A17").Value = "ok"
$ws.Range("B17").Value = 0.3359375
$ws.Range("C17").Value = 43
$ws.Range("D17").Value = 43
$ws.Range("E17").Value = 0
$ws.Range("F17").Value = 1
$ws.Range("G17").Value = $false
$ws.Range("H17").Value = 85
$ws.Range("J17").Value = "enjoy"
$ws.Range("K17").Value = 0.1666666666666667
$ws.Range("L17").Value = 31
$ws.Range("M17").Value = 31
$ws.Range("N17").Value = 1
$ws.Range("O17").Value = 0
$ws.Range("P17").Value = $false
$ws.Range("Q17").Value = 155
$ws.Range("A18").Value = "thought"
$ws.Range("B18").Value = 0.2920792079207921
$ws.Range("C18").Value = 59
$ws.Range("D18").Value = 59
$ws.Range("E18").Value = 0
$ws.Range("F18").Value = 1
$ws.Range("G18").Value = $false
$ws.Range("H18").Value = 143
$ws.Range("J18").Value = "fun"
$ws.Range("K18").Value = 0.1463628396143734
$ws.Range("L18").Value = 167
$ws.Range("M18").Value = 167
$ws.Range("N18").Value = 1
$ws.Range("O18").Value = 0
$ws.Range("P18").Value = $false
$ws.Range("Q18").Value = 974
$ws.Range("A19").Value = "size"
$ws.Range("B19").Value = 0.2474226804123711
$ws.Range("C19").Value = 48
$ws.Range("D19").Value = 48
$ws.Range("E19").Value = 0
$ws.Range("F19").Value = 1
$ws.Range("G19").Value = $false
$ws.Range("H19").Value = 146
$ws.Range("J19").Value = "easy"
$ws.Range("K19").Value = 0.08288770053475936
$ws.Range("L19").Value = 31
$ws.Range("M19").Value = 31
$ws.Range("N19").Value = 1
$ws.Range("O19").Value = 0
$ws.Range("P19").Value = $false
$ws.Range("Q19").Value = 343
$ws.Range("A20").Value = "money"
$ws.Range("B20").Value = 0.1708860759493671
$ws.Range("C20").Value = 54
$ws.Range("D20").Value = 54
$ws.Range("E20").Value = 0
$ws.Range("F20").Value = 1
$ws.Range("G20").Value = $false
$ws.Range("H20").Value = 262
$ws.Range("J20").Value = "game"
$ws.Range("K20").Value = 0.07082521117608837
$ws.Range("L20").Value = 109
$ws.Range("M20").Value = 111
$ws.Range("N20").Value = 0.98
$ws.Range("O20").Value = 0.02000000000000002
$ws.Range("P20").Value = $true
$ws.Range("Q20").Value = 1430
$ws.Range("A21").Value = "would"
$ws.Range("B21").Value = 0.1575037147102526
$ws.Range("C21").Value = 106
$ws.Range("D21").Value = 107
$ws.Range("E21").Value = 0.01
$ws.Range("F21").Value = 0.99
$ws.Range("G21").Value = $true
$ws.Range("H21").Value = 567
$ws.Range("J21").Value = "play"
$ws.Range("K21").Value = 0.04127829560585886
$ws.Range("L21").Value = 31
$ws.Range("M21").Value = 32
$ws.Range("N21").Value = 0.97
$ws.Range("O21").Value = 0.03000000000000003
$ws.Range("P21").Value = $true
$ws.Range("Q21").Value = 720
$ws.Range("A22").Value = "hard"
$ws.Range("B22").Value = 0.15
$ws.Range("C22").Value = 30
$ws.Range("D22").Value = 30
$ws.Range("E22").Value = 0
$ws.Range("F22").Value = 1
$ws.Range("G22").Value = $false
$ws.Range("H22").Value = 170
$ws.Range("A23").Value = "item"
$ws.Range("B23").Value = 0.1485507246376812
$ws.Range("C23").Value = 41
$ws.Range("D23").Value = 41
$ws.Range("E23").Value = 0
$ws.Range("F23").Value = 1
$ws.Range("G23").Value = $false
$ws.Range("H23").Value = 235
$ws.Range("A24").Value = "price"
$ws.Range("B24").Value = 0.146551724137931
$ws.Range("C24").Value = 51
$ws.Range("D24").Value = 51
$ws.Range("E24").Value = 0
$ws.Range("F24").Value = 1
$ws.Range("G24").Value = $false
$ws.Range("H24").Value = 297
$ws.Range("A25").Value = "better"
$ws.Range("B25").Value = 0.1448598130841121
$ws.Range("C25").Value = 31
$ws.Range("D25").Value = 31
$ws.Range("E25").Value = 0
$ws.Range("F25").Value = 1
$ws.Range("G25").Value = $false
$ws.Range("H25").Value = 183
$ws.Range("A26").Value = "work"
$ws.Range("B26").Value = 0.1396825396825397
$ws.Range("C26").Value = 44
$ws.Range("D26").Value = 45
$ws.Range("E26").Value = 0.02
$ws.Range("F26").Value = 0.98
$ws.Range("G26").Value = $true
$ws.Range("H26").Value = 271
$ws.Range("A27").Value = "product"
$ws.Range("B27").Value = 0.1145374449339207
$ws.Range("C27").Value = 52
$ws.Range("D27").Value = 52
$ws.Range("E27").Value = 0
$ws.Range("F27").Value = 1
$ws.Range("G27").Value = $false
$ws.Range("H27").Value = 402
$ws.Range("A28").Value = "use"
$ws.Range("B28").Value = 0.0821917808219178
$ws.Range("C28").Value = 30
$ws.Range("D28").Value = 30
$ws.Range("E28").Value = 0
$ws.Range("F28").Value = 1
$ws.Range("G28").Value = $false
$ws.Range("H28").Value = 335
$ws.Range("A29").Value = "little"
$ws.Range("B29").Value = 0.0779510022271715
$ws.Range("C29").Value = 35
$ws.Range("D29").Value = 35
$ws.Range("E29").Value = 0
$ws.Range("F29").Value = 1
$ws.Range("G29").Value = $false
$ws.Range("H29").Value = 414
$ws.Range("A30").Value = "like"
$ws.Range("B30").Value = 0.07107438016528926
$ws.Range("C30").Value = 43
$ws.Range("D30").Value = 46
$ws.Range("E30").Value = 0.07
$ws.Range("F30").Value = 0.9299999999999999
$ws.Range("G30").Value = $true
$ws.Range("H30").Value = 562
$ws.Range("A31").Value = "one"
$ws.Range("B31").Value = 0.04580152671755725
$ws.Range("C31").Value = 36
$ws.Range("D31").Value = 44
$ws.Range("E31").Value = 0.18
$ws.Range("F31").Value = 0.8200000000000001
$ws.Range("G31").Value = $true
$ws.Range("H31").Value = 750
